# Todo.docx: collapse the "Pink ghost" / "Blue ghost" / "Orange ghost"
# bullet points into a single bullet announcing that ghosts stay hidden
# until the "Player One" text disappears, while preserving the
# "_GoBack" bookmark that originally lived inside the "Blue ghost" item.

$d = $word.ActiveDocument

# Locate the three ghost bullets by content (robust to any incidental
# paragraph-index shift) instead of trusting fixed indices.
$pinkIdx = 0
$blueIdx = 0
$orangeIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Pink ghost*")   { $pinkIdx = $i }
    if ($t -like "Blue ghost*")   { $blueIdx = $i }
    if ($t -like "Orange ghost*") { $orangeIdx = $i }
}

# --- Step 1: drop the "Orange ghost ..." paragraph entirely --------------
# (it is the paragraph right after the "Blue ghost ..." one)
$pBlue   = $d.Paragraphs.Item($blueIdx)
$pOrange = $d.Paragraphs.Item($orangeIdx)
$d.Range($pBlue.Range.End, $pOrange.Range.End).Delete()

# --- Step 2: merge "Pink ghost ..." into "Blue ghost ..." ----------------
# Delete the "Pink ghost ..." paragraph's text together with its paragraph
# mark so the following paragraph (which still owns the bookmark) becomes
# the merged bullet. This keeps the bookmark run-adjacency intact instead
# of recreating the bookmark (recreating it exactly on a paragraph-mark
# boundary is unreliable); we just retarget its surrounding text below.
$pPink = $d.Paragraphs.Item($pinkIdx)
$pBlue2 = $d.Paragraphs.Item($blueIdx)
$d.Range($pPink.Range.Start, $pBlue2.Range.Start).Delete()

# --- Step 3: swap in the new sentence before the bookmark -----------------
$d.Content.Find.Execute(
    "Blue ghost should move to center of box and then up ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Make sure ghost are not shown until " + [char]8220 + "Player One" + [char]8221 + " text disappears.",
    2)

# --- Step 4: drop the remaining trailing text after the bookmark ---------
$d.Content.Find.Execute(
    "and out. Make this ghost move like the others.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)
